$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 581.5
$ws.Range("I2").Value = 372.5
$ws.Range("K2").Value = 372.5
$ws.Range("M2").Value = -259.5
$ws.Range("H19").Value = 1051.8235
$ws.Range("I19").Value = 1130.1
$ws.Range("K19").Value = 1130.1
$ws.Range("M19").Value = -955.0999999999999
$ws.Range("H38").Value = 1439.5454
$ws.Range("I38").Value = 1439.5454
$ws.Range("K38").Value = 4318.6362
$ws.Range("M38").Value = -3946.6362
$ws.Range("H43").Value = 30539.4
$ws.Range("I43").Value = 699
$ws.Range("K43").Value = 699
$ws.Range("M43").Value = -630
$ws.Range("H87").Value = 24213.84
$ws.Range("J87").Value = 24213.84
$ws.Range("L87").Value = 24213.84
$ws.Range("N87").Value = -26709.84
$ws.Range("H88").Value = 4787
$ws.Range("I88").Value = 8499.5
$ws.Range("J88").Value = 1074.5
$ws.Range("K88").Value = 8499.5
$ws.Range("L88").Value = 1074.5
$ws.Range("M88").Value = -8093.5
$ws.Range("N88").Value = -1886.5
$ws.Range("H90").Value = 24213.84
$ws.Range("J90").Value = 24213.84
$ws.Range("L90").Value = 72641.52
$ws.Range("N90").Value = -85121.52
$ws.Range("H91").Value = 4787
$ws.Range("I91").Value = 8499.5
$ws.Range("J91").Value = 1074.5
$ws.Range("K91").Value = 8499.5
$ws.Range("L91").Value = 1074.5
$ws.Range("M91").Value = -7095.5
$ws.Range("N91").Value = -3882.5
$ws.Range("H96").Value = 1233.3334
$ws.Range("I96").Value = 850
$ws.Range("K96").Value = 2550
$ws.Range("M96").Value = -1177
$ws.Range("H111").Value = 2845.32
$ws.Range("I111").Value = 2026
$ws.Range("K111").Value = 6078
$ws.Range("M111").Value = -3011
$ws.Range("H132").Value = 21005.154
$ws.Range("I132").Value = 1642.4546
$ws.Range("K132").Value = 4927.3638
$ws.Range("M132").Value = -2397.3638

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7928.0786
$ws.Range("I32").Value = 7246.66
$ws.Range("K32").Value = 7246.66
$ws.Range("M32").Value = -6959.66
$ws.Range("H45").Value = 1606.7646
$ws.Range("I45").Value = 923.36365
$ws.Range("K45").Value = 923.36365
$ws.Range("M45").Value = -546.36365
$ws.Range("H74").Value = 469.6087
$ws.Range("I74").Value = 440.7647
$ws.Range("J74").Value = 551.3333
$ws.Range("K74").Value = 440.7647
$ws.Range("L74").Value = 551.3333
$ws.Range("M74").Value = 433.2353
$ws.Range("N74").Value = -2299.3333
$ws.Range("H77").Value = 469.6087
$ws.Range("I77").Value = 440.7647
$ws.Range("J77").Value = 551.3333
$ws.Range("K77").Value = 2203.8235
$ws.Range("L77").Value = 2756.6665
$ws.Range("M77").Value = 2164.1765
$ws.Range("N77").Value = -11492.6665
$ws.Range("H97").Value = 3832230.2
$ws.Range("I97").Value = 838.087
$ws.Range("J97").Value = 18519234
$ws.Range("K97").Value = 838.087
$ws.Range("L97").Value = 18519234
$ws.Range("M97").Value = -342.087
$ws.Range("N97").Value = -18520226

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 6694.4165
$ws.Range("J86").Value = 7148.222
$ws.Range("L86").Value = 7148.222
$ws.Range("N86").Value = -9394.222
$ws.Range("H89").Value = 6694.4165
$ws.Range("J89").Value = 7148.222
$ws.Range("L89").Value = 35741.11
$ws.Range("N89").Value = -46973.11
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").Value = ""
$ws.Range("H105").Value = 4425.16
$ws.Range("I105").Value = 4017.5
$ws.Range("J105").Value = 5149.8887
$ws.Range("K105").Value = 4017.5
$ws.Range("L105").Value = 5149.8887
$ws.Range("M105").Value = -2270.5
$ws.Range("N105").Value = -8643.8887
$ws.Range("H107").Value = 972.86664
$ws.Range("I107").Value = 966.0833
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 966.0833
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 953.9167
$ws.Range("N107").Value = -4840

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 50000
$ws.Range("J20").Value = 50000
$ws.Range("L20").Value = 50000
$ws.Range("N20").Value = -50472
$ws.Range("H22").Value = 999.6667
$ws.Range("I22").Value = 198.2
$ws.Range("K22").Value = 198.2
$ws.Range("M22").Value = 151.8
$ws.Range("H30").Value = 50000
$ws.Range("J30").Value = 50000
$ws.Range("L30").Value = 50000
$ws.Range("N30").Value = -50182
$ws.Range("H31").Value = 3400.5715
$ws.Range("J31").Value = 5296.1763
$ws.Range("L31").Value = 5296.1763
$ws.Range("N31").Value = -5886.1763
$ws.Range("H34").Value = 3400.5715
$ws.Range("J34").Value = 5296.1763
$ws.Range("L34").Value = 5296.1763
$ws.Range("N34").Value = -5700.1763
$ws.Range("H99").Value = 11226361
$ws.Range("I99").Value = 2037950.6
$ws.Range("K99").Value = 2037950.6
$ws.Range("M99").Value = -2036452.6
$ws.Range("H106").Value = 27000
$ws.Range("J106").Value = 27000
$ws.Range("L106").Value = 27000
$ws.Range("N106").Value = -29524
$ws.Range("H126").Value = 11226361
$ws.Range("I126").Value = 2037950.6
$ws.Range("K126").Value = 6113851.800000001
$ws.Range("M126").Value = -6111381.800000001
$ws.Range("H128").Value = 50000
$ws.Range("J128").Value = 50000
$ws.Range("L128").Value = 50000
$ws.Range("N128").Value = -59960
$ws.Range("H134").Value = 5100.756
$ws.Range("I134").Value = 3713.4583
$ws.Range("K134").Value = 11140.3749
$ws.Range("M134").Value = -8605.374899999999
$ws.Range("H141").Value = 113738.48
$ws.Range("J141").Value = 113738.48
$ws.Range("L141").Value = 113738.48
$ws.Range("N141").Value = -124098.48

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 203.55556
$ws.Range("I6").Value = 118.85714
$ws.Range("K6").Value = 356.57142
$ws.Range("M6").Value = -243.57142
$ws.Range("H63").Value = 4985.3
$ws.Range("J63").Value = 5000
$ws.Range("L63").Value = 15000
$ws.Range("N63").Value = -16498
$ws.Range("H66").Value = 4985.3
$ws.Range("J66").Value = 5000
$ws.Range("L66").Value = 45000
$ws.Range("N66").Value = -52488
$ws.Range("H133").Value = 5137.5
$ws.Range("I133").Value = 5137.5
$ws.Range("K133").Value = 15412.5
$ws.Range("M133").Value = -10352.5
$ws.Range("H140").Value = 2008.2727
$ws.Range("I140").Value = 1984.9048
$ws.Range("J140").Value = 2499
$ws.Range("K140").Value = 5954.7144
$ws.Range("L140").Value = 7497
$ws.Range("M140").Value = -774.7143999999998
$ws.Range("N140").Value = -17857

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 426.6
$ws.Range("I107").Value = 309.42856
$ws.Range("K107").Value = 309.42856
$ws.Range("M107").Value = 1610.57144
$ws.Range("H126").Value = 5233.0835
$ws.Range("I126").Value = 4332.3335
$ws.Range("K126").Value = 12997.0005
$ws.Range("M126").Value = -10527.0005

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 4999
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").Value = ""
$ws.Range("H16").Value = 1107.1666
$ws.Range("I16").Value = 328.6
$ws.Range("J16").Value = 5000
$ws.Range("K16").Value = 328.6
$ws.Range("L16").Value = 5000
$ws.Range("M16").Value = -158.6
$ws.Range("N16").Value = -5340
$ws.Range("H48").Value = 25666.334
$ws.Range("J48").Value = 25666.334
$ws.Range("L48").Value = 25666.334
$ws.Range("N48").Value = -26988.334
$ws.Range("H132").Value = 5233.0713
$ws.Range("I132").Value = 2802
$ws.Range("J132").Value = 6205.5
$ws.Range("K132").Value = 8406
$ws.Range("L132").Value = 18616.5
$ws.Range("M132").Value = -5876
$ws.Range("N132").Value = -23676.5
$ws.Range("H136").Value = 4251.364
$ws.Range("I136").Value = 2404.524
$ws.Range("K136").Value = 7213.572
$ws.Range("M136").Value = -4663.572

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 770
$ws.Range("I7").Value = 912.5
$ws.Range("J7").Value = 200
$ws.Range("K7").Value = 912.5
$ws.Range("L7").Value = 200
$ws.Range("M7").Value = -799.5
$ws.Range("N7").Value = -426
$ws.Range("H9").Value = 10200
$ws.Range("J9").Value = 10200
$ws.Range("L9").Value = 10200
$ws.Range("N9").Value = -10480
$ws.Range("H13").Value = 6476.75
$ws.Range("I13").Value = 12250.5
$ws.Range("J13").Value = 703
$ws.Range("K13").Value = 12250.5
$ws.Range("L13").Value = 703
$ws.Range("M13").Value = -12110.5
$ws.Range("N13").Value = -983
$ws.Range("H81").Value = 23819472
$ws.Range("I81").Value = 14088.667
$ws.Range("K81").Value = 28177.334
$ws.Range("M81").Value = -27116.334
$ws.Range("H84").Value = 23819472
$ws.Range("I84").Value = 14088.667
$ws.Range("K84").Value = 140886.67
$ws.Range("M84").Value = -135582.67
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").Value = ""
$ws.Range("H122").Value = 2872.282
$ws.Range("I122").Value = 2320.6667
$ws.Range("K122").Value = 6962.000100000001
$ws.Range("M122").Value = -4512.000100000001
